$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("C37").Value = "Age at diagnosis of CVD [years]"
$ws.Range("C39").Value = "Age at diagnosis of angina pectoris [years]"
$ws.Range("C41").Value = "Age at diagnosis of myocardial infarction [years]"
$ws.Range("C43").Value = "Age at diagnosis of stroke [years]"
$ws.Range("C45").Value = "Age at diagnosis of cerebral infarction (ischaemic stroke) [years]"
$ws.Range("C47").Value = "Age at diagnosis of haemorrhagic stroke [years]"
$ws.Range("C49").Value = "Age at diagnosis of essential hypertension [years]"
$ws.Range("C51").Value = "Age at diagnosis of heart failure [years]"
$ws.Range("C53").Value = "Age at diagnosis of diabetes mellitus type 2 [years]"
$ws.Range("C56").Value = "Age at diagnosis of cancer [years]"
$ws.Range("C58").Value = "Age at time of death [years]"
